$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 167
$ws.Range("F5").Value = 3907
$ws.Range("F6").Value = 224
$ws.Range("F7").Value = 2574
$ws.Range("F9").Value = 3182
$ws.Range("F11").Value = 2332
$ws.Range("F14").Value = 329
$ws.Range("F15").Value = 466
$ws.Range("F16").Value = 16
$ws.Range("F18").Value = 218
$ws.Range("F20").Value = 311
$ws.Range("F21").Value = 432
$ws.Range("F22").Value = 670
$ws.Range("F24").Value = 48
$ws.Range("F26").Value = 1307
$ws.Range("F27").Value = 138
$ws.Range("F28").Value = 157
$ws.Range("F29").Value = 36
$ws.Range("F31").Value = 67
$ws.Range("F32").Value = 4373
$ws.Range("F33").Value = 4200
$ws.Range("F34").Value = 84
$ws.Range("F35").Value = 134
$ws.Range("F36").Value = 65
$ws.Range("F37").Value = 16
$ws.Range("F38").Value = 1149
$ws.Range("F39").Value = 8
$ws.Range("F40").Value = 486
$ws.Range("F42").Value = 1323
$ws.Range("F43").Value = 180
$ws.Range("F44").Value = 134
$ws.Range("F47").Value = 66

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 6

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 155
$ws.Range("F4").Value = 2327

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 155
$ws.Range("F8").Value = 167
$ws.Range("F10").Value = 3907
$ws.Range("F11").Value = 224
$ws.Range("F12").Value = 2574
$ws.Range("F14").Value = 3182
$ws.Range("F17").Value = 2332
$ws.Range("F20").Value = 329
$ws.Range("F21").Value = 16
$ws.Range("F24").Value = 311
$ws.Range("F25").Value = 432
$ws.Range("F26").Value = 670
$ws.Range("F28").Value = 48
$ws.Range("F29").Value = 1307
$ws.Range("F30").Value = 157
$ws.Range("F32").Value = 67
$ws.Range("F34").Value = 4373
$ws.Range("F35").Value = 4200
$ws.Range("F36").Value = 84
$ws.Range("F37").Value = 16
$ws.Range("F38").Value = 1149
$ws.Range("F39").Value = 8
$ws.Range("F45").Value = 1323
$ws.Range("F46").Value = 180
